# "Reset project to start."
# The ShapeWorks Studio project state (stored on the "studio" worksheet as
# key/value rows) is reset back to its initial values: the notes field is
# cleared out (re-saved as a fresh empty Qt rich-text blob), the active tool
# goes back to "data", the view goes back to "Original", and the zoom level
# is reset to "2".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studio")

$notes = @"
<!DOCTYPE HTML PUBLIC "-//W3C//DTD HTML 4.0//EN" "http://www.w3.org/TR/REC-html40/strict.dtd">
<html><head><meta name="qrichtext" content="1" /><style type="text/css">
p, li { white-space: pre-wrap; }
</style></head><body style=" font-family:'.AppleSystemUIFont'; font-size:13pt; font-weight:400; font-style:normal;">
<p style="-qt-paragraph-type:empty; margin-top:0px; margin-bottom:0px; margin-left:0px; margin-right:0px; -qt-block-indent:0; text-indent:0px; font-family:'Noto Sans'; font-size:11pt;"><br /></p></body></html>
"@

$ws.Range("B4").Value = $notes
# Writing the long wrapped-text note triggers an auto row-height bump; put
# row 4 back to its natural (non-custom) height so only the cell value
# changes, matching the source edit.
$ws.Rows.Item(4).AutoFit()

$ws.Range("B5").Value = "data"
$ws.Range("B6").Value = "Original"

# zoom_state's new value ("2") looks numeric, so a plain Range.Value
# assignment would store it as a number -- but the workbook keeps it as text
# (same as every other key/value row on this sheet). The "project" sheet's
# version cell already holds the text "2", so copy that cell's value across;
# this preserves the text type without touching number formats/styles.
$wb.Worksheets.Item("project").Range("B2").Copy($ws.Range("B7"))
